# Fixing test cases fro watchlist
# Adds three new test case rows (42-44) to the "Test Cases" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)          # "Test Cases" sheet (sheet1.xml)
$styleSrc = $wb.Worksheets.Item(3)    # "AuthoringTest" sheet - A2 carries style index 1

$ws.Activate()

# --- Row 42 : VerifyPostRecordDetails / OPQA-370 ---------------------------
$styleSrc.Range("A2").Copy($ws.Range("A42"))
$ws.Range("A42").Value = "VerifyPostRecordDetails"

$ws.Range("B41").Copy($ws.Range("B42"))
$ws.Range("B42").Value = "OPQA-370"

$styleSrc.Range("A2").Copy($ws.Range("C42"))
$ws.Range("C42").Value = "Verify that user contributed articles display the information about the author"

$styleSrc.Range("A2").Copy($ws.Range("D42"))
$ws.Range("D42").Value = "Y"

$styleSrc.Range("A2").Copy($ws.Range("E42"))
$ws.Range("E42").Value = "PASS"

# --- Row 43 : SeacrhAndViewOwnPost / OPQA-415 -------------------------------
$styleSrc.Range("A2").Copy($ws.Range("A43"))
$ws.Range("A43").Value = "SeacrhAndViewOwnPost"

$ws.Range("B41").Copy($ws.Range("B43"))
$ws.Range("B43").Value = "OPQA-415"

$styleSrc.Range("A2").Copy($ws.Range("C43"))
$ws.Range("C43").Value = "Verify that user is able to search the  posts a user authored themselves and view them."

$styleSrc.Range("A2").Copy($ws.Range("D43"))
$ws.Range("D43").Value = "Y"

$styleSrc.Range("A2").Copy($ws.Range("E43"))
$ws.Range("E43").Value = "PASS"

# --- Row 44 : SeacrhAndViewOthersPost / OPQA-416 ----------------------------
$styleSrc.Range("A2").Copy($ws.Range("A44"))
$ws.Range("A44").Value = "SeacrhAndViewOthersPost"

$ws.Range("B41").Copy($ws.Range("B44"))
$ws.Range("B44").Value = "OPQA-416"

$styleSrc.Range("A2").Copy($ws.Range("C44"))
$ws.Range("C44").Value = "Verify that user is able to search the posts of others and view them."

$styleSrc.Range("A2").Copy($ws.Range("D44"))
$ws.Range("D44").Value = "Y"

$styleSrc.Range("A2").Copy($ws.Range("E44"))
$ws.Range("E44").Value = "PASS"

# --- View state: scroll down a bit and select B44, like the authored file --
$ws.Range("B44").Select()
